$wb = $excel.ActiveWorkbook

# Sheet ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1563.8
$ws.Range("L17").Value = 4691.4
$ws.Range("J17").Value = 1563.8
$ws.Range("N17").Value = -5027.4

# Sheet ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M33").Value = -413.1429000000001
$ws.Range("H33").Value = 642.1429000000001
$ws.Range("I33").Value = 642.1429000000001
$ws.Range("K33").Value = 642.1429000000001

# Sheet ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3267
$ws.Range("L100").Value = 2866.3333
$ws.Range("J100").Value = 2866.3333
$ws.Range("N100").Value = -3948.3333

# Sheet ALC row 115
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

# Sheet ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M131").Value = 1524.75
$ws.Range("H131").Value = 1737.4
$ws.Range("I131").Value = 1171.75
$ws.Range("K131").Value = 3515.25

# Sheet ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M132").Value = -9477.8465
$ws.Range("H132").Value = 3905.8572
$ws.Range("I132").Value = 4002.6155
$ws.Range("K132").Value = 12007.8465

# Sheet ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M2").Value = -799.34784
$ws.Range("H2").Value = 890.92
$ws.Range("L2").Value = 644.5
$ws.Range("J2").Value = 644.5
$ws.Range("I2").Value = 912.34784
$ws.Range("N2").Value = -870.5
$ws.Range("K2").Value = 912.34784

# Sheet ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M61").Value = -1606.8667
$ws.Range("H61").Value = 1818.8667
$ws.Range("I61").Value = 1818.8667
$ws.Range("K61").Value = 1818.8667

# Sheet ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M110").Value = -2899.385
$ws.Range("H110").Value = 4944.385
$ws.Range("I110").Value = 4944.385
$ws.Range("K110").Value = 4944.385

# Sheet ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M116").Value = 1381.65216
$ws.Range("H116").Value = 890.92
$ws.Range("L116").Value = 644.5
$ws.Range("J116").Value = 644.5
$ws.Range("I116").Value = 912.34784
$ws.Range("N116").Value = -5232.5
$ws.Range("K116").Value = 912.34784

# Sheet ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M132").Value = -6065.500100000001
$ws.Range("H132").Value = 4316.5884
$ws.Range("I132").Value = 2865.1667
$ws.Range("K132").Value = 8595.500100000001

# Sheet ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M136").Value = -2906.6001
$ws.Range("H136").Value = 1818.8667
$ws.Range("I136").Value = 1818.8667
$ws.Range("K136").Value = 5456.6001

# Sheet BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M3").Value = -798.34784
$ws.Range("H3").Value = 890.92
$ws.Range("L3").Value = 644.5
$ws.Range("J3").Value = 644.5
$ws.Range("I3").Value = 912.34784
$ws.Range("N3").Value = -872.5
$ws.Range("K3").Value = 912.34784

# Sheet BSM row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 7531.625
$ws.Range("L80").Value = 11890
$ws.Range("J80").Value = 11890
$ws.Range("N80").Value = -13886

# Sheet BSM row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 7531.625
$ws.Range("L83").Value = 59450
$ws.Range("J83").Value = 11890
$ws.Range("N83").Value = -69434

# Sheet BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M105").Value = -2183.7693
$ws.Range("H105").Value = 4673.8076
$ws.Range("I105").Value = 3930.7693
$ws.Range("K105").Value = 3930.7693

# Sheet CRP row 3
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M3").Value = -108.5
$ws.Range("H3").Value = 481
$ws.Range("I3").Value = 221.5
$ws.Range("K3").Value = 221.5

# Sheet CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M31").Value = -2301
$ws.Range("H31").Value = 2596
$ws.Range("I31").Value = 2596
$ws.Range("K31").Value = 2596

# Sheet CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M34").Value = -2394
$ws.Range("H34").Value = 2596
$ws.Range("I34").Value = 2596
$ws.Range("K34").Value = 2596

# Sheet CRP row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M51").Value = -2197
$ws.Range("H51").Value = 2933
$ws.Range("I51").Value = 2933
$ws.Range("K51").Value = 2933

# Sheet CRP row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M61").Value = -2585
$ws.Range("H61").Value = 2933
$ws.Range("I61").Value = 2933
$ws.Range("K61").Value = 2933

# Sheet CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M86").Value = -6856.8
$ws.Range("H86").Value = 7984.0835
$ws.Range("L86").Value = 7987.143
$ws.Range("J86").Value = 7987.143
$ws.Range("I86").Value = 7979.8
$ws.Range("N86").Value = -10233.143
$ws.Range("K86").Value = 7979.8

# Sheet CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M89").Value = -34283
$ws.Range("H89").Value = 7984.0835
$ws.Range("L89").Value = 39935.715
$ws.Range("J89").Value = 7987.143
$ws.Range("I89").Value = 7979.8
$ws.Range("N89").Value = -51167.715
$ws.Range("K89").Value = 39899

# Sheet CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M99").Value = -563.75
$ws.Range("H99").Value = 2061.75
$ws.Range("L99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("I99").Value = 2061.75
$ws.Range("K99").Value = 2061.75
$ws.Range("N99").ClearContents()

# Sheet CRP row 117
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H117").Value = 59998
$ws.Range("L117").Value = 59998
$ws.Range("J117").Value = 59998
$ws.Range("N117").Value = -69176

# Sheet CRP row 120
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H120").Value = 63387.668
$ws.Range("L120").Value = 70065.2
$ws.Range("J120").Value = 70065.2
$ws.Range("N120").Value = -77323.2

# Sheet CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M122").Value = -1604.9998
$ws.Range("H122").Value = 1351.6666
$ws.Range("I122").Value = 1351.6666
$ws.Range("K122").Value = 4054.9998

# Sheet CRP row 125
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("N125").ClearContents()

# Sheet CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M126").Value = -3715.25
$ws.Range("H126").Value = 2061.75
$ws.Range("L126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("I126").Value = 2061.75
$ws.Range("N126").Value = 0
$ws.Range("K126").Value = 6185.25

# Sheet CUL row 20
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M20").Value = -787
$ws.Range("H20").Value = 338
$ws.Range("L20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("I20").Value = 338
$ws.Range("K20").Value = 1014
$ws.Range("N20").ClearContents()

# Sheet CUL row 22
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 51
$ws.Range("L22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("N22").ClearContents()

# Sheet CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M23").Value = -962
$ws.Range("H23").Value = 516.8
$ws.Range("L23").Value = 2080.5
$ws.Range("J23").Value = 693.5
$ws.Range("I23").Value = 399
$ws.Range("N23").Value = -2550.5
$ws.Range("K23").Value = 1197

# Sheet CUL row 26
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M26").Value = -3963.6
$ws.Range("H26").Value = 1307.8572
$ws.Range("L26").Value = 1189.99998
$ws.Range("J26").Value = 396.66666
$ws.Range("I26").Value = 1417.2
$ws.Range("N26").Value = -1765.99998
$ws.Range("K26").Value = 4251.6

# Sheet CUL row 27
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 51
$ws.Range("L27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("N27").ClearContents()

# Sheet CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 8617.111000000001
$ws.Range("L39").Value = 25851.333
$ws.Range("J39").Value = 8617.111000000001
$ws.Range("N39").Value = -26439.333

# Sheet CUL row 47
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M47").Value = -943
$ws.Range("H47").Value = 458
$ws.Range("I47").Value = 458
$ws.Range("K47").Value = 1374

# Sheet CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 242.16667
$ws.Range("L92").Value = 746.00001
$ws.Range("J92").Value = 248.66667
$ws.Range("N92").Value = -3242.00001

# Sheet CUL row 94
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 3000
$ws.Range("L94").Value = 9000
$ws.Range("J94").Value = 3000
$ws.Range("N94").Value = -10352

# Sheet CUL row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1057.25
$ws.Range("L98").Value = 3729
$ws.Range("J98").Value = 1243
$ws.Range("N98").Value = -6725

# Sheet GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M70").Value = -9729.5
$ws.Range("H70").Value = 9999.5
$ws.Range("I70").Value = 9999.5
$ws.Range("K70").Value = 9999.5

# Sheet GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M73").Value = -9063.5
$ws.Range("H73").Value = 9999.5
$ws.Range("I73").Value = 9999.5
$ws.Range("K73").Value = 9999.5

# Sheet LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 785
$ws.Range("L40").Value = 785
$ws.Range("J40").Value = 785
$ws.Range("N40").Value = -1057

# Sheet LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M61").Value = -1198
$ws.Range("H61").Value = 1400
$ws.Range("I61").Value = 1400
$ws.Range("K61").Value = 1400

# Sheet LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M93").Value = 373.25
$ws.Range("H93").Value = 874.75
$ws.Range("L93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("I93").Value = 874.75
$ws.Range("K93").Value = 874.75
$ws.Range("N93").ClearContents()

# Sheet LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M113").Value = 770
$ws.Range("H113").Value = 1400
$ws.Range("I113").Value = 1400
$ws.Range("K113").Value = 1400

# Sheet LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M122").Value = -8727.571599999999
$ws.Range("H122").Value = 3907.5
$ws.Range("I122").Value = 3725.8572
$ws.Range("K122").Value = 11177.5716

# Sheet WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M81").Value = -4072
$ws.Range("H81").Value = 3469.8
$ws.Range("L81").Value = 9649.5
$ws.Range("J81").Value = 4824.75
$ws.Range("I81").Value = 2566.5
$ws.Range("N81").Value = -11771.5
$ws.Range("K81").Value = 5133

# Sheet WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M84").Value = -20361
$ws.Range("H84").Value = 3469.8
$ws.Range("L84").Value = 48247.5
$ws.Range("J84").Value = 4824.75
$ws.Range("I84").Value = 2566.5
$ws.Range("N84").Value = -58855.5
$ws.Range("K84").Value = 25665

# Sheet WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M132").Value = -7055.3531
$ws.Range("H132").Value = 3106.5
$ws.Range("I132").Value = 3195.1177
$ws.Range("K132").Value = 9585.3531
